$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the start time of row 6 (Work on Alice) from 14:00 to 14:10
$ws.Range("A6").Value = 0.59027777777777779

# Add a new row 13: copy formatting from row 9 (time-formatted, centered cell)
# then give it the new values: time 14:13 and text "Test 1 Bro"
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 0.59236111111111112
$ws.Range("B13").Value = "Test 1 Bro"

# Update the selected cell to reflect where the user ended up (B16)
$ws.Range("B16").Select()
